$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ICP-MS"

$ws.Activate()
$ws.Range("C6").Select()
